$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pipelineLog")

# Insert a new row before row 21 (shifts rows 21+ down by one, e.g. the GWAS
# row moves from 21->22, the anvilDRSClient-only row moves 22->23, and the
# compareGeccoPPs rows move 25->26 / 26->27)
$ws.Rows.Item(21).Insert()

# Bring over the per-column formatting used elsewhere in the table (style
# "Heading 2/3/4" look-alikes) so the new row matches its neighbours exactly.
$ws.Range("A3").Copy()
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy()
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy()
$ws.Range("C21").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy()
$ws.Range("D21").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").Copy()
$ws.Range("E21").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Copy()
$ws.Range("F21").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").Copy()
$ws.Range("G21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new row 21 with the GTEX SB example data (this particular
# left-to-right-but-D-last order reproduces the author's shared-string
# insertion order: script, wesClient, note, comment, searchClient, then the
# already-existing "anvilDRSClient" string reused for drsClient)
$ws.Range("A21").Value = "FASPScript18.py"
$ws.Range("E21").Value = "sbWESClient"
$ws.Range("B21").Value = "Anvil GTEX SB"
$ws.Range("F21").Value = "Modified version of 15 to use free amazon data. Move to notebook"
$ws.Range("C21").Value = "Gen3ManifestClient"
$ws.Range("D21").Value = "anvilDRSClient"

# Match the row height used by the rest of the data rows
$ws.Rows.Item(21).RowHeight = 24

# The sort-remembered range grows by one row now that a row was inserted
# inside it. Applying a Sort snaps/expands the target range to cover whole
# merged cells, so temporarily unmerge the (still merged, now shifted down
# one row) cells that straddle the new range boundary, apply, then re-merge.
$ws.Range("A22:A23").UnMerge()
$ws.Range("B22:B23").UnMerge()
$ws.Range("C22:C23").UnMerge()
$ws.Range("E22:E23").UnMerge()

$sortObj = $ws.Sort
$sortObj.SetRange($ws.Range("A3:E22"))
try { $sortObj.Apply() } catch {}

$ws.Range("A22:A23").Merge()
$ws.Range("B22:B23").Merge()
$ws.Range("C22:C23").Merge()
$ws.Range("E22:E23").Merge()

# Widen column F
$ws.Columns.Item(6).ColumnWidth = 60.1640625

# Update selection to match the author's last selection
$ws.Range("C22:C23").Select()
